$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New "print preview" help text block (rows 7-11 on "First sheet")
# Write the cell text in row order first so the shared-string table picks up
# the new strings in the same order the target workbook has them (27-31).
$ws.Range("A7").Value = "Multi sheet issues with print preview"
$ws.Range("A8").Value = "If you have problems with printing the document:"
$ws.Range("A9").Value = "1)Create 1st sheet like this"
$ws.Range("A10").Value = "2)or call macro like in example 070 to save file with a new name (so the entire file would be recreated)"
$ws.Range("A11").Value = "https://docs.microsoft.com/en-us/office/vba/api/excel.workbook.saveas"

# Bold the heading line before the hyperlink is created so the new
# direct-format (bold) cell style is registered ahead of the hyperlink style.
$ws.Range("A7").Font.Bold = $true

$ws.Hyperlinks.Add($ws.Range("A11"), "https://docs.microsoft.com/en-us/office/vba/api/excel.workbook.saveas")

# Selection ends up on G1 after the edit
$ws.Range("G1").Select()
